$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy format from H1 (style index 1:
# bold, bordered, center/top aligned) then overwrite with the new text.
$ws.Range("H1:H1").Copy($ws.Range("I1"))
$ws.Range("H1:H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for new columns I and J, rows 2-23
$data = @{
    2  = @(8, 9)
    3  = @(9, 9)
    4  = @(9, 9)
    5  = @(6, 6)
    6  = @(10, 10)
    7  = @(7, 7)
    8  = @(7, 8)
    9  = @(8, 8)
    10 = @(7, 7)
    11 = @(6, 6)
    12 = @(4, 6)
    13 = @(6, 6)
    14 = @(5, 5)
    15 = @(4, 5)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(2, 3)
    19 = @(4, 6)
    20 = @(9, 9)
    21 = @(8, 8)
    22 = @(9, 9)
    23 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
